# Feito do Exercício 5 ao 7
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mark exercises 5, 6 and 7 (rows 6-8) as done (TRUE) in column C
$ws.Range("C6").Value = $true
$ws.Range("C7").Value = $true
$ws.Range("C8").Value = $true

# Move the active selection to B16 (last cell touched by the author)
$ws.Range("B16").Select() | Out-Null
